{"js": "// Update the narrative paragraph text describing the amended charges.\nconst narrativeResults = context.document.body.search(\n  \"is amended to Possession of Marijuana Drug Paraphernalia\",\n  { matchCase: true }\n);\nawait context.sync();\n\nif (narrativeResults.items.length > 0) {\n  narrativeResults.items[0].insertText(\n    \"is amended to Driving Under Suspension FTA, Fines or Child Support and Tail Lights-rear License Plate is amended to Traffic Control Device\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// Update the sentencing table: Offense / Statute-Ord. / Degree rows.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n\n  // Offense row (row 0): column 1 (DUS Ucm) and column 2 (Tail Lights).\n  table.getCell(0, 1).value =\n    \"DUS Ucm - AMENDED to Driving Under Suspension FTA, Fines or Child Support\";\n  table.getCell(0, 2).value =\n    \"Tail Lights-rear License Plate - AMENDED to Traffic Control Device\";\n\n  // Statute/Ord. row (row 1).\n  table.getCell(1, 1).value = \"4510.111 \";\n  table.getCell(1, 2).value = \"4511.12\";\n\n  // Degree row (row 2).\n  table.getCell(2, 1).value = \"Unclassified Misdemeanor\";\n  table.getCell(2, 2).value = \"Minor Misdemeanor\";\n\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the narrative paragraph describing the amended charges.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Execute(\n    \"is amended to Possession of Marijuana Drug Paraphernalia\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"is amended to Driving Under Suspension FTA, Fines or Child Support and Tail Lights-rear License Plate is amended to Traffic Control Device\",\n    2\n) | Out-Null\n\n# Update the sentencing table: Offense / Statute-Ord. / Degree rows.\n$table = $d.Tables.Item(1)\n\n# Offense row (row 1): column 2 (DUS Ucm) and column 3 (Tail Lights).\n$table.Cell(1, 2).Range.Text = \"DUS Ucm - AMENDED to Driving Under Suspension FTA, Fines or Child Support\"\n$table.Cell(1, 3).Range.Text = \"Tail Lights-rear License Plate - AMENDED to Traffic Control Device\"\n\n# Statute/Ord. row (row 2).\n$table.Cell(2, 2).Range.Text = \"4510.111 \"\n$table.Cell(2, 3).Range.Text = \"4511.12\"\n\n# Degree row (row 3).\n$table.Cell(3, 2).Range.Text = \"Unclassified Misdemeanor\"\n$table.Cell(3, 3).Range.Text = \"Minor Misdemeanor\"\n"}
